$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.539.86'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.907.01'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").Value = "'" + '1.009'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'" + '338.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.19%  '
$ws.Range("D6").Value = "'" + '1.008'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -1.20%  '
$ws.Range("D8").Value = "'" + '0.4000'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("E9").Value = '  -2.44%  '
$ws.Range("D10").Value = "'" + '0.9905'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("D11").Value = "'" + '23.30'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("D12").Value = '1.915.63'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("E13").Value = '  -2.66%  '
$ws.Range("D14").Value = "'" + '7.109'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").Value = "'" + '89.12'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("D16").Value = "'" + '0.06839'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = "'" + '0.00001022'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.90%  '
$ws.Range("D19").Value = "'" + '17.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = "'" + '1.007'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '29.529.46'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = "'" + '5.507'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = "'" + '11.59'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").Value = '2.147.12'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").Value = "'" + '156.93'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").Value = "'" + '19.59'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").Value = "'" + '2.054'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("D30").Value = "'" + '119.21'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = "'" + '0.9942'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.49%  '
$ws.Range("D32").Value = "'" + '0.09527'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("D33").Value = "'" + '5.486'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("D34").Value = "'" + '3.547'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").Value = "'" + '0.06457'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.68%  '
$ws.Range("D37").Value = "'" + '0.02238'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("D38").Value = "'" + '1.193'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("D39").Value = "'" + '0.5798'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  -3.14%  '
$ws.Range("D41").Value = "'" + '7.746'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.29%  '
$ws.Range("D42").Value = "'" + '0.1823'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").Value = "'" + '2.451'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = "'" + '0.07428'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'" + '12.13'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.14%  '
$ws.Range("D47").Value = "'" + '0.5476'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("D48").Value = "'" + '1.946'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("D49").Value = "'" + '116.12'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").Value = "'" + '2.375'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = "'" + '71.24'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.62%  '
